$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "26.951.51"
Set-TextValue "E2" "  +1.14%  "
Set-TextValue "D3" "1.846.72"
Set-TextValue "E3" "  +1.04%  "
Set-TextValue "E4" "  +0.41%  "
Set-TextValue "E5" "  +0.34%  "
Set-TextValue "D6" "309.75"
Set-TextValue "E6" "  +0.08%  "
Set-TextValue "D7" "0.4778"
Set-TextValue "E7" "  +2.39%  "
Set-TextValue "E8" "  +2.14%  "
Set-TextValue "D9" "0.07232"
Set-TextValue "E9" "  +1.30%  "
Set-TextValue "D10" "0.9293"
Set-TextValue "E10" "  +2.81%  "
Set-TextValue "D12" "0.07728"
Set-TextValue "E12" "  +0.53%  "
Set-TextValue "D13" "1.822.72"
Set-TextValue "E13" "  -0.89%  "
Set-TextValue "D14" "5.337"
Set-TextValue "E14" "  +1.33%  "
Set-TextValue "D15" "6.430"
Set-TextValue "E15" "  +0.82%  "
Set-TextValue "D16" "88.68"
Set-TextValue "E17" "  +0.45%  "
Set-TextValue "D18" "0.000008644"
Set-TextValue "E18" "  +1.02%  "
Set-TextValue "E19" "  +0.35%  "
Set-TextValue "D20" "26.984.54"
Set-TextValue "E20" "  +1.10%  "
Set-TextValue "D21" "14.45"
Set-TextValue "E21" "  +1.70%  "
Set-TextValue "D22" "5.058"
Set-TextValue "E22" "  +0.65%  "
Set-TextValue "D24" "1.924"
Set-TextValue "E24" "  +0.70%  "
Set-TextValue "D25" "152.77"
Set-TextValue "E25" "  -0.12%  "
Set-TextValue "E26" "  +1.48%  "
Set-TextValue "D27" "2.001"
Set-TextValue "E27" "  +0.16%  "
Set-TextValue "D28" "114.28"
Set-TextValue "E28" "  +0.35%  "
Set-TextValue "D29" "4.964"
Set-TextValue "E29" "  +2.07%  "
Set-TextValue "D30" "0.08891"
Set-TextValue "E30" "  +0.75%  "
Set-TextValue "D31" "3.318"
Set-TextValue "E31" "  +5.46%  "
Set-TextValue "E32" "  +0.64%  "
Set-TextValue "D33" "0.7425"
Set-TextValue "E33" "  +0.70%  "
Set-TextValue "D34" "4.496"
Set-TextValue "E34" "  +1.37%  "
Set-TextValue "D35" "2.752"
Set-TextValue "E35" "  -3.96%  "
Set-TextValue "D36" "1.119"
Set-TextValue "E36" "  +3.65%  "
Set-TextValue "D37" "0.01955"
Set-TextValue "E37" "  +1.09%  "
Set-TextValue "D38" "0.05268"
Set-TextValue "E38" "  +2.04%  "
Set-TextValue "D39" "2.980"
Set-TextValue "E39" "  +1.40%  "
Set-TextValue "D40" "0.5216"
Set-TextValue "E40" "  +2.86%  "
Set-TextValue "D41" "6.987"
Set-TextValue "E41" "  +1.55%  "
Set-TextValue "D42" "0.1511"
Set-TextValue "E42" "  +0.81%  "
Set-TextValue "D43" "8.210"
Set-TextValue "E43" "  +1.83%  "
Set-TextValue "E44" "  +6.26%  "
Set-TextValue "E45" "  +1.90%  "
Set-TextValue "E46" "  +0.38%  "
Set-TextValue "D47" "101.68"
Set-TextValue "E47" "  +3.20%  "
Set-TextValue "D48" "1.605"
Set-TextValue "D49" "65.68"
Set-TextValue "E49" "  +2.71%  "
Set-TextValue "D50" "0.06060"
Set-TextValue "E50" "  +0.62%  "
Set-TextValue "D51" "0.8875"
Set-TextValue "E51" "  +4.04%  "
